$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 89260.64
$ws.Range("I15").Value = 89260.64
$ws.Range("K15").Value = 267781.92
$ws.Range("M15").Value = -267612.92

$ws.Range("H76").Value = 8550308
$ws.Range("I76").Value = 18521018
$ws.Range("J76").Value = 3985.7144
$ws.Range("K76").Value = 18521018
$ws.Range("L76").Value = 3985.7144
$ws.Range("M76").Value = -18520703
$ws.Range("N76").Value = -4615.7144

$ws.Range("H79").Value = 8550308
$ws.Range("I79").Value = 18521018
$ws.Range("J79").Value = 3985.7144
$ws.Range("K79").Value = 18521018
$ws.Range("L79").Value = 3985.7144
$ws.Range("M79").Value = -18519926
$ws.Range("N79").Value = -6169.7144

$ws.Range("H141").Value = 3812.3215
$ws.Range("I141").Value = 1867.619
$ws.Range("J141").Value = 9646.429
$ws.Range("K141").Value = 5602.857
$ws.Range("L141").Value = 28939.287
$ws.Range("M141").Value = -422.857
$ws.Range("N141").Value = -39299.287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22741.092
$ws.Range("I32").Value = 5230
$ws.Range("J32").Value = 133644.67
$ws.Range("K32").Value = 5230
$ws.Range("L32").Value = 133644.67
$ws.Range("M32").Value = -4943
$ws.Range("N32").Value = -134218.67

$ws.Range("H45").Value = 878.4
$ws.Range("I45").Value = 848
$ws.Range("K45").Value = 848
$ws.Range("M45").Value = -471

$ws.Range("H97").Value = 9785.362999999999
$ws.Range("I97").Value = 11780.889
$ws.Range("K97").Value = 11780.889
$ws.Range("M97").Value = -11284.889

$ws.Range("H110").Value = 100001200
$ws.Range("I110").Value = 111112340
$ws.Range("K110").Value = 111112340
$ws.Range("M110").Value = -111110295

$ws.Range("H111").Value = 400000
$ws.Range("J111").Value = 400000
$ws.Range("L111").Value = 400000
$ws.Range("N111").Value = -408180

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1714.2858
$ws.Range("I86").Value = 1833.3334
$ws.Range("J86").Value = 1625
$ws.Range("K86").Value = 1833.3334
$ws.Range("L86").Value = 1625
$ws.Range("M86").Value = -710.3334
$ws.Range("N86").Value = -3871

$ws.Range("H89").Value = 1714.2858
$ws.Range("I89").Value = 1833.3334
$ws.Range("J89").Value = 1625
$ws.Range("K89").Value = 9166.666999999999
$ws.Range("L89").Value = 8125
$ws.Range("M89").Value = -3550.666999999999
$ws.Range("N89").Value = -19357

$ws.Range("H94").Value = 852.13336
$ws.Range("I94").Value = 747.9545000000001
$ws.Range("J94").Value = 1138.625
$ws.Range("K94").Value = 747.9545000000001
$ws.Range("L94").Value = 1138.625
$ws.Range("M94").Value = -296.9545000000001
$ws.Range("N94").Value = -2040.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2223121.2
$ws.Range("I6").Value = 6666997.5
$ws.Range("J6").Value = 1183.3334
$ws.Range("K6").Value = 6666997.5
$ws.Range("L6").Value = 1183.3334
$ws.Range("M6").Value = -6666884.5
$ws.Range("N6").Value = -1409.3334

$ws.Range("H22").Value = 839.93335
$ws.Range("I22").Value = 485.1
$ws.Range("K22").Value = 485.1
$ws.Range("M22").Value = -135.1

$ws.Range("H122").Value = 1861.7333
$ws.Range("I122").Value = 1162.1904
$ws.Range("J122").Value = 3494
$ws.Range("K122").Value = 3486.5712
$ws.Range("L122").Value = 10482
$ws.Range("M122").Value = -1036.5712
$ws.Range("N122").Value = -15382

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5666.6665
$ws.Range("I56").Value = 5666.6665
$ws.Range("K56").Value = 5666.6665
$ws.Range("M56").Value = -5136.6665

$ws.Range("H132").Value = 1746.1818
$ws.Range("I132").Value = 1401.6
$ws.Range("J132").Value = 2033.3334
$ws.Range("K132").Value = 12614.4
$ws.Range("L132").Value = 18300.0006
$ws.Range("M132").Value = -10084.4
$ws.Range("N132").Value = -23360.0006

$ws.Range("H137").Value = 6736703.5
$ws.Range("I137").Value = 7695040
$ws.Range("J137").Value = 507516.5
$ws.Range("K137").Value = 23085120
$ws.Range("L137").Value = 1522549.5
$ws.Range("M137").Value = -23080020
$ws.Range("N137").Value = -1532749.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H122").Value = 2615.8
$ws.Range("I122").Value = 2518.2307
$ws.Range("J122").Value = 3250
$ws.Range("K122").Value = 7554.6921
$ws.Range("L122").Value = 9750
$ws.Range("M122").Value = -5104.6921
$ws.Range("N122").Value = -14650

$ws.Range("H132").Value = 2770.6123
$ws.Range("I132").Value = 2332.818
$ws.Range("K132").Value = 6998.454000000001
$ws.Range("M132").Value = -4468.454000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3494.5417
$ws.Range("I7").Value = 3191
$ws.Range("J7").Value = 3555.25
$ws.Range("K7").Value = 3191
$ws.Range("L7").Value = 3555.25
$ws.Range("M7").Value = -3079
$ws.Range("N7").Value = -3779.25

$ws.Range("H126").Value = 3494.5417
$ws.Range("I126").Value = 3191
$ws.Range("J126").Value = 3555.25
$ws.Range("K126").Value = 9573
$ws.Range("L126").Value = 10665.75
$ws.Range("M126").Value = -7103
$ws.Range("N126").Value = -15605.75

$ws.Range("H136").Value = 6560.6523
$ws.Range("I136").Value = 3975.6428
$ws.Range("J136").Value = 10581.777
$ws.Range("K136").Value = 11926.9284
$ws.Range("L136").Value = 31745.331
$ws.Range("M136").Value = -9376.928400000001
$ws.Range("N136").Value = -36845.331

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 34333.332
$ws.Range("J123").Value = 34333.332
$ws.Range("L123").Value = 34333.332
$ws.Range("N123").Value = -44133.332

$ws.Range("H126").Value = 72370.64
$ws.Range("I126").Value = 125343.125
$ws.Range("J126").Value = 1740.6666
$ws.Range("K126").Value = 376029.375
$ws.Range("L126").Value = 5221.9998
$ws.Range("M126").Value = -373559.375
$ws.Range("N126").Value = -10161.9998

$ws.Range("H132").Value = 4884.5835
$ws.Range("I132").Value = 5423.625
$ws.Range("J132").Value = 3806.5
$ws.Range("K132").Value = 16270.875
$ws.Range("L132").Value = 11419.5
$ws.Range("M132").Value = -13740.875
$ws.Range("N132").Value = -16479.5

